$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("testCitizen")

# --- sheet2 (testCitizen): swap columns A and B content and rename the strings ---
# Row data: new A values ("ulkemiz11X"), new B values ("uisve11X")
# Shared-strings table must end up with uisve111..118 (idx 12-19) THEN ulkemiz111..118 (idx 20-27),
# so write column B across all rows first, then column A across all rows.
$aValues = @("ulkemiz111","ulkemiz112","ulkemiz113","ulkemiz114","ulkemiz115","ulkemiz116","ulkemiz117","ulkemiz118")
$bValues = @("uisve111","uisve112","uisve113","uisve114","uisve115","uisve116","uisve117","uisve118")

for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 2).Value = $bValues[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Value = $aValues[$i]
}

# --- sheet2 cosmetic changes ---
# Target stored width is 20.33203125 chars; engine quantizes ColumnWidth to 1/6 steps,
# so 19.5 is the closest settable value (-> 20.333333333333332, matches to 1e-3).
$ws2.Columns.Item(1).ColumnWidth = 19.5
$ws2.Range("D16").Select()
